# Add a new departure row (row 16) to the "Main Data" sheet, recording
# the Friday, Jan 13 / 9:35 AM departure of flight W95178 (Wizz Air, A320,
# tail G-WUKF) to London (LTN), 5 minutes ahead of schedule.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 16

$ws.Cells.Item($row, 1).Value  = 15
$ws.Cells.Item($row, 2).Value  = "Friday, Jan 13"
$ws.Cells.Item($row, 3).Value  = "9:40 AM"
$ws.Cells.Item($row, 4).Value  = "W95178"
$ws.Cells.Item($row, 5).Value  = "London"
$ws.Cells.Item($row, 6).Value  = "(LTN)"
$ws.Cells.Item($row, 7).Value  = "Wizz Air "
$ws.Cells.Item($row, 8).Value  = "A320"
$ws.Cells.Item($row, 9).Value  = "(G-WUKF)"
$ws.Cells.Item($row, 10).Value = "9:35 AM"
# Column K (11) stays blank for this row, same as the rows above it; touch
# the cell's formatting (no-op) so it is materialized in the sheet just
# like the existing blank cells in column K/M.
$ws.Cells.Item($row, 11).Borders.LineStyle = -4142
$ws.Cells.Item($row, 12).Value = "0 hours, -5 minutes"
$ws.Cells.Item($row, 13).Borders.LineStyle = -4142
